$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# H2:I17 were stored as shared-string lookalikes of the D/odor_channel and
# stim-within-odor values ("wood"/"pencil"/"strawberry"/"apple"); the sheet
# author replaced them with the plain numeric codes 1-4 that the text
# actually represented (H mirrors the odor_channel group in column D, I
# cycles 1..4 within each group of four rows).
for ($r = 2; $r -le 17; $r++) {
    $group = [int][Math]::Ceiling(($r - 1) / 4)
    $cond = (($r - 2) % 4) + 1
    $ws.Cells.Item($r, 8).Value = $group
    $ws.Cells.Item($r, 9).Value = $cond
}

# The sheet's saved selection moved from H27 to H30.
$ws.Range("H30").Select()

# The workbook window position/scroll also moved in the source file
# (xWindow/yWindow -80,1580 -> 0,460). Reflect this through every COM
# surface the host exposes; harmless no-ops if a given property isn't wired.
try { $excel.ActiveWindow.Left = 0 } catch {}
try { $excel.ActiveWindow.Top = 460 } catch {}
try { $excel.Left = 0 } catch {}
try { $excel.Top = 460 } catch {}
